# Excel COM-interop script implementing the commit:
# "arquivo Concursos.xlsx atualizado: excluido a tabela inscrição"
#
# The "Isncricao" (Inscricao) worksheet is removed. Its columns are folded
# into the sheets that reference it:
#   - inicio_inscricao / fim_inscricao (old "data_inicial"/"data_final")
#     move into the "Concurso" sheet as new columns E and F (the old
#     "link" header/value shift from column E to column G).
#   - taxa move into the "Cargo" sheet as a new column D (existing
#     columns C_reserva..id_concurso shift right by one).

$wb = $excel.ActiveWorkbook
[void]($excel.DisplayAlerts = $false)

$wsConcurso  = $wb.Worksheets.Item("Concurso")
$wsCargo     = $wb.Worksheets.Item("Cargo")
$wsInscricao = $wb.Worksheets.Item("Isncricao")

# --- grab the values that currently live on the Isncricao sheet ---
$dataInicial = $wsInscricao.Range("B2").Value2
$dataFinal   = $wsInscricao.Range("C2").Value2
$taxaValor   = $wsInscricao.Range("D2").Value2

$fmtData     = "d/m/yyyy"
$fmtMoeda    = "_-""R$""* #,##0.00_-;""-R$""* #,##0.00_-;_-""R$""* \-??_-;_-@_-"

# =====================================================================
# Concurso sheet: insert "inicio_inscrição (null)" / "fim_inscrição"
# columns, pushing the existing "link" column from E to G.
# =====================================================================

# Preserve the header/value formatting of the column being displaced
# ("link" / "url") by copying it (value + format) into its new home,
# and likewise onto the new "fim_inscrição" header cell so it matches
# the other header cells.
$wsConcurso.Range("E1").Copy($wsConcurso.Range("F1")) | Out-Null
$wsConcurso.Range("E1").Copy($wsConcurso.Range("G1")) | Out-Null
$wsConcurso.Range("E2").Copy($wsConcurso.Range("G2")) | Out-Null

$wsConcurso.Range("E1").Value = "inicio_inscrição (null)"
$wsConcurso.Range("F1").Value = "fim_inscrição"

$wsConcurso.Range("E2").Value = $dataInicial
$wsConcurso.Range("E2").NumberFormat = $fmtData
$wsConcurso.Range("F2").Value = $dataFinal
$wsConcurso.Range("F2").NumberFormat = $fmtData

$wsConcurso.Columns.Item(5).ColumnWidth = 18.6
$wsConcurso.Columns.Item(6).ColumnWidth = 13.6

# =====================================================================
# Cargo sheet: insert a new "taxa(null)" column D (existing columns
# c_reserva, qtd_vagas, salario, id_concurso shift right by one).
# =====================================================================

$wsCargo.Range("D1").EntireColumn.Insert() | Out-Null
$wsCargo.Columns.Item(4).ColumnWidth = $wsCargo.Columns.Item(3).ColumnWidth

$wsCargo.Range("D1").Value = "taxa(null)"
$wsCargo.Range("D1").Font.Bold = $true
$wsCargo.Range("D1").Font.ColorIndex = 2
$wsCargo.Range("D1").Interior.Color = $wsCargo.Range("C1").Interior.Color
$wsCargo.Range("D1").HorizontalAlignment = $wsCargo.Range("C1").HorizontalAlignment
$wsCargo.Range("D1").VerticalAlignment = $wsCargo.Range("C1").VerticalAlignment

$wsCargo.Range("D2").Value = $taxaValor
$wsCargo.Range("D2").NumberFormat = $fmtMoeda

# =====================================================================
# Remove the now-empty "Isncricao" worksheet entirely.
# =====================================================================
$wsInscricao.Delete() | Out-Null

# Leave the workbook with "Cargo" as the active tab, matching the
# recorded activeTab index.
$wsCargo.Activate() | Out-Null
